# Apply strikethrough formatting to every bullet item under "Control
# Validation" EXCEPT the two "Date Controls" bullets, matching the
# target revision described by the diff.
#
# Paragraph map (1-based, via $d.Paragraphs):
#   1  NSCC Hackathon 2025 / Control Validation   <- untouched (title)
#   2  Employee Name
#   3  Letters only, one space between words
#   4  Format <letter-string><space><letter-string> only
#   5  Max 40 characters
#   6  W#
#   7  Must be in format "wXXXXXXX"
#   8  Max 8 characters
#   9  Leading 'w' and numbers only
#   10 Fund/Dept/Program/Acct/Project
#   11 Alphanumeric only
#   12 Max 40(?) characters
#   13 Pay start
#   14 Selected date must be a Sunday
#   15 Date Controls                              <- untouched
#   16 All date controls must fall within ...      <- untouched
#   17 Hours Worked
#   18 Float  >= 0 only
#   19 Max 3 characters (XX.X)
#   20 Other Information
#   21 Disallow semi-colons
#   22 Max 128 characters
#   23 Notes/Comments
#   24 Disallow semi-colons
#   25 Max 255 characters

$d = $word.ActiveDocument

# --- Step 1: bulk strikethrough ------------------------------------------
# Applying Font.StrikeThrough on the paragraph's own Range (which, in
# Word's object model, includes the trailing paragraph mark) sets strike
# both on the paragraph mark's run properties (w:pPr/w:rPr/w:strike) and
# on every run of text in the paragraph (w:r/w:rPr/w:strike) - exactly the
# pattern seen throughout the diff.

for ($i = 2; $i -le 14; $i++) {
    $d.Paragraphs.Item($i).Range.Font.StrikeThrough = 1
}

for ($i = 17; $i -le 25; $i++) {
    $d.Paragraphs.Item($i).Range.Font.StrikeThrough = 1
}

# --- Step 2: fix up the two paragraphs whose runs/proofing marks were ----
# also reshuffled in the target revision (this happens in the real Word
# UI as a side effect of the background grammar/spell checker re-scanning
# text after it is edited). We reproduce that exact markup with
# Range.InsertXML, replacing only the paragraph's text content (not its
# trailing paragraph mark, so list numbering/pPr stays intact).

# Paragraph 18: "Float  >= 0 only"
$pFloat = $d.Paragraphs.Item(18)
$rFloat = $d.Range($pFloat.Range.Start, $pFloat.Range.End - 1)
$xmlFloat = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Float </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve"> &gt;</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:strike/></w:rPr><w:t>= 0 only</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rFloat.InsertXML($xmlFloat)

# Paragraph 25: "Max 255 characters"
$pMax255 = $d.Paragraphs.Item(25)
$rMax255 = $d.Range($pMax255.Range.Start, $pMax255.Range.End - 1)
$xmlMax255 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:strike/></w:rPr><w:t xml:space="preserve">Max 255 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:strike/></w:rPr><w:t>characters</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rMax255.InsertXML($xmlMax255)
